# Vietnamese (vi-VN) localization of the built-in placeholder / layout /
# notes-master shape names and the instructional placeholder text that
# ships in PowerPoint's default template, matching the upstream
# document-templates commit that produced vi-VN/new.pptx.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Common instructional-text translations, reused across every shape
# that still carries the stock English placeholder copy.
# ---------------------------------------------------------------------
$TitleStyleText    = "Bấm để sửa kiểu tiêu đề Bản cái"
$SubtitleStyleText = "Bấm & sửa kiểu phụ đề của Bản chính"
$BodyLvl1Text      = "Bấm để sửa kiểu văn bản Bản cái"
$Lvl2Text          = "Mức hai"
$Lvl3Text          = "Mức ba"
$Lvl4Text          = "Mức bốn"
$Lvl5Text          = "Mức năm"
$PictureIconText   = "Nhấp vào biểu tượng để thêm hình ảnh"

function Set-BodyLevels {
    param($shape)
    $tr = $shape.TextFrame.TextRange
    $tr.Paragraphs(1,1).Text = $BodyLvl1Text
    $tr.Paragraphs(2,1).Text = $Lvl2Text
    $tr.Paragraphs(3,1).Text = $Lvl3Text
    $tr.Paragraphs(4,1).Text = $Lvl4Text
    $tr.Paragraphs(5,1).Text = $Lvl5Text
}

# ---------------------------------------------------------------------
# 1. ppt/slides/slide1.xml — only the shape names change, no body text.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).Name = "Tiêu đề 1"
$s1.Shapes.Item(2).Name = "Tiêu đề phụ 2"

# ---------------------------------------------------------------------
# 2. ppt/notesSlides/notesSlide1.xml — only shape names change.
# ---------------------------------------------------------------------
$np = $s1.NotesPage
$np.Shapes.Item(1).Name = "Slide Hình ảnh giữ chỗ 1"
$np.Shapes.Item(2).Name = "Ghi chú Chỗ giữ chỗ 2"
$np.Shapes.Item(3).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 5"

# ---------------------------------------------------------------------
# 3. ppt/notesMasters/notesMaster1.xml — only shape names change.
# ---------------------------------------------------------------------
$nm = $p.NotesMaster
$nm.Shapes.Item(1).Name = "Tiêu đề giữ chỗ 1"
$nm.Shapes.Item(2).Name = "Chỗ dành sẵn cho Ngày tháng 2"
$nm.Shapes.Item(3).Name = "Slide Hình ảnh giữ chỗ 3"
$nm.Shapes.Item(4).Name = "Ghi chú Chỗ giữ chỗ 4"
$nm.Shapes.Item(5).Name = "Chỗ dành sẵn cho Chân trang 5"
$nm.Shapes.Item(6).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 6"

# ---------------------------------------------------------------------
# 4. ppt/slideMasters/slideMaster1.xml — names + instructional text.
# ---------------------------------------------------------------------
$sm = $p.SlideMaster
$sm.Shapes.Item(1).Name = "Chỗ dành sẵn cho Tiêu đề 1"
$sm.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1).Text = $TitleStyleText

$sm.Shapes.Item(2).Name = "Chỗ dành sẵn cho Văn bản 2"
Set-BodyLevels $sm.Shapes.Item(2)

$sm.Shapes.Item(3).Name = "Chỗ dành sẵn cho Ngày tháng 3"
$sm.Shapes.Item(4).Name = "Chỗ dành sẵn cho Chân trang 4"
$sm.Shapes.Item(5).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 5"

# ---------------------------------------------------------------------
# 5. ppt/slideLayouts/slideLayout{1..11}.xml — cSld name, shape names
#    and instructional text.
# ---------------------------------------------------------------------
$layouts = $sm.CustomLayouts

# ---- Layout 1: Title Slide -> Tiêu đề Bản chiếu -----------------------
$lay = $layouts.Item(1)
$lay.Name = "Tiêu đề Bản chiếu"
$lay.Shapes.Item(1).Name = "Tiêu đề 1"
$lay.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1).Text = $TitleStyleText
$lay.Shapes.Item(2).Name = "Tiêu đề phụ 2"
$lay.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1,1).Text = $SubtitleStyleText
$lay.Shapes.Item(3).Name = "Chỗ dành sẵn cho Ngày tháng 3"
$lay.Shapes.Item(4).Name = "Chỗ dành sẵn cho Chân trang 4"
$lay.Shapes.Item(5).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 5"

# ---- Layout 2: Title and Content -> Tiêu đề và Nội dung ---------------
$lay = $layouts.Item(2)
$lay.Name = "Tiêu đề và Nội dung"
$lay.Shapes.Item(1).Name = "Tiêu đề 1"
$lay.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1).Text = $TitleStyleText
$lay.Shapes.Item(2).Name = "Chỗ dành sẵn cho Nội dung 2"
Set-BodyLevels $lay.Shapes.Item(2)
$lay.Shapes.Item(3).Name = "Chỗ dành sẵn cho Ngày tháng 3"
$lay.Shapes.Item(4).Name = "Chỗ dành sẵn cho Chân trang 4"
$lay.Shapes.Item(5).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 5"

# ---- Layout 3: Section Header -> Đầu trang của Phần -------------------
$lay = $layouts.Item(3)
$lay.Name = "Đầu trang của Phần"
$lay.Shapes.Item(1).Name = "Tiêu đề 1"
$lay.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1).Text = $TitleStyleText
$lay.Shapes.Item(2).Name = "Chỗ dành sẵn cho Văn bản 2"
$lay.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1,1).Text = $BodyLvl1Text
$lay.Shapes.Item(3).Name = "Chỗ dành sẵn cho Ngày tháng 3"
$lay.Shapes.Item(4).Name = "Chỗ dành sẵn cho Chân trang 4"
$lay.Shapes.Item(5).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 5"

# ---- Layout 4: Two Content -> Hai Nội dung -----------------------------
$lay = $layouts.Item(4)
$lay.Name = "Hai Nội dung"
$lay.Shapes.Item(1).Name = "Tiêu đề 1"
$lay.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1).Text = $TitleStyleText
$lay.Shapes.Item(2).Name = "Chỗ dành sẵn cho Nội dung 2"
Set-BodyLevels $lay.Shapes.Item(2)
$lay.Shapes.Item(3).Name = "Chỗ dành sẵn cho Nội dung 3"
Set-BodyLevels $lay.Shapes.Item(3)
$lay.Shapes.Item(4).Name = "Chỗ dành sẵn cho Ngày tháng 4"
$lay.Shapes.Item(5).Name = "Chỗ dành sẵn cho Chân trang 5"
$lay.Shapes.Item(6).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 6"

# ---- Layout 5: Comparison -> Phép so sánh ------------------------------
$lay = $layouts.Item(5)
$lay.Name = "Phép so sánh"
$lay.Shapes.Item(1).Name = "Tiêu đề 1"
$lay.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1).Text = $TitleStyleText
$lay.Shapes.Item(2).Name = "Chỗ dành sẵn cho Văn bản 2"
$lay.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1,1).Text = $BodyLvl1Text
$lay.Shapes.Item(3).Name = "Chỗ dành sẵn cho Nội dung 3"
Set-BodyLevels $lay.Shapes.Item(3)
$lay.Shapes.Item(4).Name = "Chỗ dành sẵn cho Văn bản 4"
$lay.Shapes.Item(4).TextFrame.TextRange.Paragraphs(1,1).Text = $BodyLvl1Text
$lay.Shapes.Item(5).Name = "Chỗ dành sẵn cho Nội dung 5"
Set-BodyLevels $lay.Shapes.Item(5)
$lay.Shapes.Item(6).Name = "Chỗ dành sẵn cho Ngày tháng 6"
$lay.Shapes.Item(7).Name = "Chỗ dành sẵn cho Chân trang 7"
$lay.Shapes.Item(8).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 8"

# ---- Layout 6: Title Only -> Chỉ Tiêu đề -------------------------------
$lay = $layouts.Item(6)
$lay.Name = "Chỉ Tiêu đề"
$lay.Shapes.Item(1).Name = "Tiêu đề 1"
$lay.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1).Text = $TitleStyleText
$lay.Shapes.Item(2).Name = "Chỗ dành sẵn cho Ngày tháng 2"
$lay.Shapes.Item(3).Name = "Chỗ dành sẵn cho Chân trang 3"
$lay.Shapes.Item(4).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 4"

# ---- Layout 7: Blank -> Trống ------------------------------------------
$lay = $layouts.Item(7)
$lay.Name = "Trống"
$lay.Shapes.Item(1).Name = "Chỗ dành sẵn cho Ngày tháng 1"
$lay.Shapes.Item(2).Name = "Chỗ dành sẵn cho Chân trang 2"
$lay.Shapes.Item(3).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 3"

# ---- Layout 8: Content with Caption -> Nội dung với Phụ đề ------------
$lay = $layouts.Item(8)
$lay.Name = "Nội dung với Phụ đề"
$lay.Shapes.Item(1).Name = "Tiêu đề 1"
$lay.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1).Text = $TitleStyleText
$lay.Shapes.Item(2).Name = "Chỗ dành sẵn cho Nội dung 2"
Set-BodyLevels $lay.Shapes.Item(2)
$lay.Shapes.Item(3).Name = "Chỗ dành sẵn cho Văn bản 3"
$lay.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1,1).Text = $BodyLvl1Text
$lay.Shapes.Item(4).Name = "Chỗ dành sẵn cho Ngày tháng 4"
$lay.Shapes.Item(5).Name = "Chỗ dành sẵn cho Chân trang 5"
$lay.Shapes.Item(6).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 6"

# ---- Layout 9: Picture with Caption -> Ảnh với Phụ đề ------------------
$lay = $layouts.Item(9)
$lay.Name = "Ảnh với Phụ đề"
$lay.Shapes.Item(1).Name = "Tiêu đề 1"
$lay.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1).Text = $TitleStyleText
$lay.Shapes.Item(2).Name = "Chỗ dành sẵn cho Hình ảnh 2"
$lay.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1,1).Text = $PictureIconText
$lay.Shapes.Item(3).Name = "Chỗ dành sẵn cho Văn bản 3"
$lay.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1,1).Text = $BodyLvl1Text
$lay.Shapes.Item(4).Name = "Chỗ dành sẵn cho Ngày tháng 4"
$lay.Shapes.Item(5).Name = "Chỗ dành sẵn cho Chân trang 5"
$lay.Shapes.Item(6).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 6"

# ---- Layout 10: Title and Vertical Text -> Tiêu đề và Văn bản Dọc -----
$lay = $layouts.Item(10)
$lay.Name = "Tiêu đề và Văn bản Dọc"
$lay.Shapes.Item(1).Name = "Tiêu đề 1"
$lay.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1).Text = $TitleStyleText
$lay.Shapes.Item(2).Name = "Chỗ dành sẵn cho Văn bản Dọc 2"
Set-BodyLevels $lay.Shapes.Item(2)
$lay.Shapes.Item(3).Name = "Chỗ dành sẵn cho Ngày tháng 3"
$lay.Shapes.Item(4).Name = "Chỗ dành sẵn cho Chân trang 4"
$lay.Shapes.Item(5).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 5"

# ---- Layout 11: Vertical Title and Text -> Tiêu đề Dọc và Văn bản -----
$lay = $layouts.Item(11)
$lay.Name = "Tiêu đề Dọc và Văn bản"
$lay.Shapes.Item(1).Name = "Tiêu đề Dọc 1"
$lay.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1).Text = $TitleStyleText
$lay.Shapes.Item(2).Name = "Chỗ dành sẵn cho Văn bản Dọc 2"
Set-BodyLevels $lay.Shapes.Item(2)
$lay.Shapes.Item(3).Name = "Chỗ dành sẵn cho Ngày tháng 3"
$lay.Shapes.Item(4).Name = "Chỗ dành sẵn cho Chân trang 4"
$lay.Shapes.Item(5).Name = "Chỗ dành sẵn cho Số hiệu Bản chiếu 5"
